$wb = $excel.ActiveWorkbook

# Add a new worksheet. Excel inserts new sheets before the currently
# active sheet, so this becomes the first tab (sheetId 2, rId1) while
# the original "Tabelle1" becomes the second tab (rId2).
$qs = $wb.Worksheets.Add()
$qs.Name = "QS"

# Populate the new "QS" sheet with its data (row 1: 1, 0.3, 0.3).
$qs.Range("A1").Value = 1
$qs.Range("B1").Value = 0.3
$qs.Range("C1").Value = 0.3

# Match the saved selection/active cell on the new active "QS" sheet.
[void]$qs.Range("I4").Select()
